# Applies a cyclic permutation of the "species data" block between a set
# of rows on the active sheet. For each group of rows, row N's species
# data (Id, Taxonsorteringsordning, TaxonId, Artnamn, Vetenskapligt namn,
# Auktor, Ost, Nord, Publik kommentar) is replaced by the data that used
# to live on a different row in the same group, per the target mapping.
# All other columns (Rödlistade, Lokalnamn, Noggrannhet, Län, Kommun,
# Provins, Socken, Startdatum, Slutdatum, Ej återfunnen, ...) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as one unit ("species data" for a row).
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

# The "Auktor" column (H) is always emitted by the export even when blank
# (an empty-but-present cell), unlike "Publik kommentar" (AC) which is
# omitted entirely when there is no comment. Track that so a move that
# clears H still leaves a present-but-empty cell behind instead of
# deleting it outright.
$alwaysPresentCols = @("H")

# destRow -> sourceRow (source row's BEFORE-edit data becomes dest row's
# AFTER-edit data). Rows not listed are left untouched.
$mapping = @{
    11 = 12
    12 = 11
    13 = 14
    14 = 13
    39 = 41
    40 = 39
    41 = 40
    42 = 45
    43 = 42
    44 = 43
    45 = 44
}

# Snapshot the "before" values for every row referenced, since the
# mapping reads from rows that will themselves be overwritten. Use
# `.Value2` (not `.Value`) - this runtime's `.Value` getter does not
# resolve to the underlying scalar when read into a PS variable.
$snapshot = @{}
foreach ($row in ($mapping.Values | Select-Object -Unique)) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rng = $ws.Range("$col$row")
        $rowData[$col] = $rng.Value2
    }
    $snapshot[[int]$row] = $rowData
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[[int]$srcRow]
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$destRow")
        $val = $srcData[$col]
        if ($val -eq $null) {
            $cell.Value2 = ""
            if ($alwaysPresentCols -contains $col) {
                # Force the cell to stay materialized (present, empty)
                # instead of being dropped from the sheet entirely.
                $cell.Style = "Normal"
            }
        } else {
            $cell.Value2 = $val
        }
    }
}
